$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F, shifting the existing
# "District"/"Name" column (old F) to column G.
$ws.Columns("F:F").Insert()

# New header for the inserted column.
$ws.Range("F2").Value = "Address"

# Populate the Address column with the school/location portion
# of each teacher's combined Names cell (district dropped).
$ws.Range("F4").Value = 'G H S Arakeri Yelburga'
$ws.Range("F5").Value = 'S G R S High School kanakagiri'
$ws.Range("F6").Value = 'G H S MallapurGangavathi'
$ws.Range("F7").Value = 'G H S NiloorAfzalpur'
$ws.Range("F8").Value = 'G H S BuragapalliSedam'
$ws.Range("F9").Value = 'G J C Bhagyanagar'
$ws.Range("F10").Value = 'B G H S Tavaragera Kustagi'
$ws.Range("F11").Value = 'G H S MannurAfzalpur'
$ws.Range("F12").Value = 'G H S Kanagadda'
$ws.Range("F13").Value = 'G H S HosurAfzalpur'
$ws.Range("F14").Value = 'G H S StationGanagapurAfzalpur'
$ws.Range("F15").Value = 'G H S YadagaSedam'
$ws.Range("F16").Value = 'G P U College HanamanalKushtagi'
$ws.Range("F17").Value = 'Adarsha Vidyalaya Taluvagera'
$ws.Range("F18").Value = 'G H S GhattargaAfzalpur'
$ws.Range("F19").Value = 'G H S GanadhalYelaburga'
$ws.Range("F20").Value = 'G H S KinnisultanAland'
$ws.Range("F21").Value = 'G H S AloorJewargi'
$ws.Range("F22").Value = 'G H S ShadipurChincholi'
$ws.Range("F23").Value = 'B Shyamsunder Memorial High School Rajapur'
$ws.Range("F25").Value = 'G H S BalabattiJewargi'
$ws.Range("F26").Value = 'G H S KoradakeraKushtagi'
$ws.Range("F27").Value = 'G H S (RMSA) BaragurGangavathi'
$ws.Range("F28").Value = 'G H S Taralakatti Yalaburga'
$ws.Range("F29").Value = 'A V (RMSA) ChennurJewargi'
$ws.Range("F30").Value = 'G H S RudnoorChincholi'
$ws.Range("F31").Value = 'Sri Sevaniketan High School Melkunda(B)'
$ws.Range("F32").Value = 'H R Sarojamma Comp. Jr. Coll. Gangavathi'
$ws.Range("F33").Value = 'G G H S MudholSedam'
$ws.Range("F34").Value = 'G H S SangolgiAland'
$ws.Range("F35").Value = 'G H S MatakiAland'
$ws.Range("F36").Value = 'G H S(RMSA) ShivoorAfzalpura'
$ws.Range("F37").Value = 'G H S TalakeriYelburga'
$ws.Range("F38").Value = 'G H S KaramudiYelaburga'
$ws.Range("F39").Value = 'G H S UdagiSedam'
$ws.Range("F40").Value = 'G H S MudenoorKushtagi'
